$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed sample-set values (mean calculation added upstream to the
# dataset pipeline, shifting every per-threshold accuracy figure in this
# sheet). Row 1 keeps its "HK_R_acc_LT" label; rows 2-49 get the refreshed
# figures.
$ws.Range("A1").Value = "HK_R_acc_LT"

$ws.Range("A2").Value = 80.850091407678249
$ws.Range("A3").Value = 83.912248628884825
$ws.Range("A4").Value = 84.232175502742237
$ws.Range("A5").Value = 87.979890310786104
$ws.Range("A6").Value = 87.979890310786104
$ws.Range("A7").Value = 87.568555758683729
$ws.Range("A8").Value = 78.702010968921385
$ws.Range("A9").Value = 80.484460694698356
$ws.Range("A10").Value = 79.478976234003653
$ws.Range("A11").Value = 78.47349177330895
$ws.Range("A12").Value = 74.177330895795251
$ws.Range("A13").Value = 79.707495429616088
$ws.Range("A14").Value = 77.787934186471659
$ws.Range("A15").Value = 78.107861060329071
$ws.Range("A16").Value = 80.393053016453393
$ws.Range("A17").Value = 77.60511882998172
$ws.Range("A18").Value = 79.75319926873857
$ws.Range("A19").Value = 85.009140767824505
$ws.Range("A20").Value = 87.934186471663622
$ws.Range("A21").Value = 87.934186471663622
$ws.Range("A22").Value = 87.751371115173669
$ws.Range("A23").Value = 79.433272394881172
$ws.Range("A24").Value = 82.038391224862877
$ws.Range("A25").Value = 81.444241316270578
$ws.Range("A26").Value = 81.581352833638022
$ws.Range("A27").Value = 80.758683729433272
$ws.Range("A28").Value = 81.627056672760517
$ws.Range("A29").Value = 81.215722120658143
$ws.Range("A30").Value = 80.21023765996344
$ws.Range("A31").Value = 80.118829981718463
$ws.Range("A32").Value = 89.396709323583181
$ws.Range("A33").Value = 92.413162705667276
$ws.Range("A34").Value = 92.001828153564901
$ws.Range("A35").Value = 83.135283363802557
$ws.Range("A36").Value = 87.979890310786104
$ws.Range("A37").Value = 68.00731261425959
$ws.Range("A38").Value = 84.597806215722116
$ws.Range("A39").Value = 80.804387568555754
$ws.Range("A40").Value = 79.524680073126149
$ws.Range("A41").Value = 79.387568555758676
$ws.Range("A42").Value = 79.478976234003653
$ws.Range("A43").Value = 79.661791590493607
$ws.Range("A44").Value = 79.616087751371111
$ws.Range("A45").Value = 82.952468007312618
$ws.Range("A46").Value = 84.643510054844612
$ws.Range("A47").Value = 78.93053016453382
$ws.Range("A48").Value = 77.787934186471659
$ws.Range("A49").Value = 80.393053016453393
